$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 97 (existing rows 97-152 shift down to 99-154).
$ws.Rows.Item(97).Resize(2).Insert()

# New row 97: Especial quality record for date 2021-11-10 (serial 44510).
$ws.Range("A97").Value = 3
$ws.Range("B97").Value = "Femacal de La Calera"
$ws.Range("C97").Value = "Coquimbo"
$ws.Range("D97").Value = 44510
$ws.Range("E97").Value = 5
$ws.Range("F97").Value = "Fruta"
$ws.Range("G97").Value = 100101
$ws.Range("H97").Value = "Berries"
$ws.Range("I97").Value = 100112025
$ws.Range("J97").Value = "Frutilla"
$ws.Range("K97").Value = "Sin especificar"
$ws.Range("L97").Value = "Especial"
$ws.Range("M97").Value = 85
$ws.Range("N97").Value = 6000
$ws.Range("O97").Value = 6000
$ws.Range("P97").Value = 6000
$ws.Range("Q97").Value = "$/bandeja 7 kilos"
$ws.Range("R97").Value = "Provincia de Melipilla"
$ws.Range("S97").Value = 857
$ws.Range("T97").Value = 7

# New row 98: Segunda quality record for the same date.
$ws.Range("A98").Value = 3
$ws.Range("B98").Value = "Femacal de La Calera"
$ws.Range("C98").Value = "Coquimbo"
$ws.Range("D98").Value = 44510
$ws.Range("E98").Value = 5
$ws.Range("F98").Value = "Fruta"
$ws.Range("G98").Value = 100101
$ws.Range("H98").Value = "Berries"
$ws.Range("I98").Value = 100112025
$ws.Range("J98").Value = "Frutilla"
$ws.Range("K98").Value = "Sin especificar"
$ws.Range("L98").Value = "Segunda"
$ws.Range("M98").Value = 75
$ws.Range("N98").Value = 4000
$ws.Range("O98").Value = 4000
$ws.Range("P98").Value = 4000
$ws.Range("Q98").Value = "$/bandeja 7 kilos"
$ws.Range("R98").Value = "Provincia de Melipilla"
$ws.Range("S98").Value = 571
$ws.Range("T98").Value = 7
